# Add I0 and IF columns (I and J) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the formatting of H1 (bold header style) onto the new
# header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-35.
$values = @(
    @(9, 9),
    @(5, 7),
    @(5, 6),
    @(8, 9),
    @(3, 4),
    @(4, 5),
    @(4, 6),
    @(1, 3),
    @(4, 6),
    @(1, 1),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(10, 10),
    @(5, 5),
    @(6, 6),
    @(6, 7),
    @(7, 8),
    @(5, 6),
    @(8, 8),
    @(6, 6),
    @(10, 10),
    @(11, 11),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(8, 8),
    @(5, 5),
    @(7, 7),
    @(9, 9),
    @(5, 5),
    @(8, 8),
    @(4, 4),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
